# Generate Report for Handoff
#
# The localization-status report is regenerated and the "Overview" sheet's
# "Latest HO Xliff Generate Date" column for the 601907bc-...md row (row 7)
# picks up a fresh handoff timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("G7").Value = "2016-08-25 04:40:34"
